# Uruguay Primera Division - atualizacao de bases (19-04-2024 23:27)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write one "match row" worth of cells at a given row index.
# $cols is an ordered hashtable-like array of (ColumnLetter, Value) pairs.
# ---------------------------------------------------------------------------
function Set-RowCells {
    param($Row, $Values)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}
function Get-RowSnapshot {
    param($Row)
    $snap = [ordered]@{}
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")) {
        $cell = $ws.Range("$col$Row")
        $snap[$col] = $cell.Value2
    }
    return $snap
}

# ---------------------------------------------------------------------------
# 1) Rows 117-119: three fixtures get re-sorted (cyclic rotation). The "id"
#    in column A stays fixed per row position; every other column moves:
#      new row117 <- old row118
#      new row118 <- old row119
#      new row119 <- old row117
# ---------------------------------------------------------------------------

$row117 = [ordered]@{
    B = 7013886; F = "Racing Club de Montevideo"; G = "Cerro"; H = 0; I = 1; J = "A";
    K = 2.25; L = 3.1; M = 3.25; N = 2.25; O = 2.875; P = 3.5; Q = -0.25; R = 1.95; S = 1.9;
    T = 2; U = 1.925; V = 1.925; W = -1; X = -1; Y = 2.5; Z = -1; AA = 0.8999999999999999; AB = -1; AC = 0.925
}
$row118 = [ordered]@{
    B = 7013702; F = "Defensor Sporting"; G = "Danubio"; H = 0; I = 2; J = "A";
    K = 1.8; L = 3.6; M = 4.2; N = 1.8; O = 3.6; P = 4.2; Q = -0.75; R = 2.05; S = 1.8;
    T = 2.25; U = 1.85; V = 2; W = -1; X = -1; Y = 3.2; Z = -1; AA = 0.8; AB = -0.5; AC = 0.5
}
$row119 = [ordered]@{
    B = 7013885; F = "La Luz"; G = "Atletico Fenix Montevideo"; H = 0; I = 2; J = "A";
    K = 3; L = 3; M = 2.4; N = 2.9; O = 2.75; P = 2.6; Q = 0; R = 2.025; S = 1.825;
    T = 2; U = 2.025; V = 1.825; W = -1; X = -1; Y = 1.6; Z = -1; AA = 0.825; AB = 0; AC = 0
}

Set-RowCells 117 $row117
Set-RowCells 118 $row118
Set-RowCells 119 $row119

# ---------------------------------------------------------------------------
# 2) Two brand new fixtures are inserted right before the old row 185,
#    pushing the old rows 185/186/187 down to 187/188/189. A further new
#    fixture is appended at the very end, as the new row 190.
#
#    We achieve the "insert" purely through value writes (bottom-up) so we
#    never disturb existing styles/shared strings via a COM Insert call:
#      new190 <- brand new fixture (8081250)
#      new189 <- old187 (8081249)
#      new188 <- old186 (8081144)
#      new187 <- old185 (8081162)
#      new186 <- brand new fixture (8081435)
#      new185 <- brand new fixture (8081163)
# ---------------------------------------------------------------------------

# Snapshot the old content of rows 185-187 (as plain values) before overwriting.
$old185 = Get-RowSnapshot 185
$old186 = Get-RowSnapshot 186
$old187 = Get-RowSnapshot 187

# -- new row 190 (brand-new fixture, appended at the end) -------------------
# Copy formatting (styles) from row 187 (an existing fully-formatted data
# row) down onto the three freshly-used rows 188/189/190 first.
$ws.Range("A187:AC187").Copy() | Out-Null
$ws.Range("A188:AC190").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$row190 = [ordered]@{
    A = 188; B = 8081250; C = "Uruguay Primera División"; D = "Uruguay Apertura";
    E = 45402.85416666666; F = "Deportivo Maldonado"; G = "Miramar Misiones";
    K = 2.2; L = 3.3; M = 3.3; N = 2.1; O = 3.3; P = 3.6; Q = -0.25; R = 1.85; S = 2;
    T = 2.25; U = 2; V = 1.85; W = 0; X = 0; Y = 0; Z = 0; AA = 0
}
Set-RowCells 190 $row190

# -- new row 189 <- old187 ----------------------------------------------------
Set-RowCells 189 $old187

# -- new row 188 <- old186 ----------------------------------------------------
Set-RowCells 188 $old186

# -- new row 187 <- old185 ----------------------------------------------------
Set-RowCells 187 $old185

# -- new row 186 (brand-new fixture, inserted) -------------------------------
# row 185/186 need H/I/J cells too (the two new fixtures carry a final
# result, unlike the neighbouring still-unplayed fixtures) and style copies.
$ws.Range("A183:AC183").Copy() | Out-Null
$ws.Range("A185:AC186").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$row186 = [ordered]@{
    A = 184; B = 8081435; C = "Uruguay Primera División"; D = "Uruguay Apertura";
    E = 45401.64583333334; F = "Liverpool Montevideo"; G = "Defensor Sporting";
    H = 1; I = 1; J = "D";
    K = 2.375; L = 3.3; M = 3; N = 3.1; O = 3.5; P = 2.2; Q = 0.25; R = 1.925; S = 1.925;
    T = 2.5; U = 2.05; V = 1.8; W = -1; X = 2.5; Y = -1; Z = 0.4625; AA = -0.5; AB = -1; AC = 0.8
}
Set-RowCells 186 $row186

# -- new row 185 (brand-new fixture, inserted) -------------------------------
$row185 = [ordered]@{
    A = 183; B = 8081163; C = "Uruguay Primera División"; D = "Uruguay Apertura";
    E = 45401.5625; F = "Racing Club de Montevideo"; G = "Cerro";
    H = 1; I = 1; J = "H";
    K = 2.3; L = 3.2; M = 3.2; N = 2.2; O = 3.2; P = 3.5; Q = -0.25; R = 1.85; S = 2;
    T = 2.25; U = 2.05; V = 1.8; W = -1; X = 2.2; Y = -1; Z = -0.5; AA = 0.5; AB = -0.5; AC = 0.4
}
Set-RowCells 185 $row185
